# Re-order the product rows (2-7) in a cyclic fashion so that each row's
# full record moves to its new position, and then sync the "Item No."
# (column A) and "Mfr Catalog No." (column B) columns: whichever column is
# missing gets filled from the other; if both ended up empty, leave both
# empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current block (rows 2-7, columns A-R) before touching anything.
$src = $ws.Range("A2:R7").Value2

$rowCount = 6
$colCount = 18

# destMap[destRowIndex] = sourceRowIndex (both 1-based, relative to the block)
#   dest row 2 (index 1) <- source row 7 (index 6)
#   dest row 3 (index 2) <- source row 5 (index 4)
#   dest row 4 (index 3) <- source row 3 (index 2)
#   dest row 5 (index 4) <- source row 4 (index 3)
#   dest row 6 (index 5) <- source row 6 (index 5)   (unchanged)
#   dest row 7 (index 6) <- source row 2 (index 1)
$destMap = @{ 1 = 6; 2 = 4; 3 = 2; 4 = 3; 5 = 5; 6 = 1 }

$out = New-Object 'object[,]' $rowCount, $colCount

for ($destIdx = 1; $destIdx -le $rowCount; $destIdx++) {
    $srcIdx = $destMap[$destIdx]
    for ($col = 1; $col -le $colCount; $col++) {
        $out[$destIdx - 1, $col - 1] = $src[$srcIdx, $col]
    }

    # Column A is index 1, column B is index 2 within the block.
    $aVal = $out[$destIdx - 1, 0]
    $bVal = $out[$destIdx - 1, 1]

    $aEmpty = ($null -eq $aVal) -or ($aVal -eq "")
    $bEmpty = ($null -eq $bVal) -or ($bVal -eq "") -or ($bVal -eq "N/A")

    if (-not $aEmpty -and $bEmpty) {
        # Item No. present, Mfr Catalog No. missing -> copy it across.
        $out[$destIdx - 1, 1] = $aVal
    } elseif ($aEmpty -and $bEmpty) {
        # Neither is meaningfully populated -> clear both.
        $out[$destIdx - 1, 0] = ""
        $out[$destIdx - 1, 1] = ""
    }
}

# Column N ("# List Price") holds currency-looking text such as "$56.82".
# Assigning that through .Value2 normally auto-coerces it into a number, so
# mark just that column as Text first to keep it as literal string content,
# then restore General formatting afterwards (the underlying cell already
# holds a string by then, so the format flip back does not re-coerce it).
$priceRange = $ws.Range("N2:N7")
$priceRange.NumberFormat = "@"

$ws.Range("A2:R7").Value2 = $out

$priceRange.NumberFormat = "General"
